# Add three new BOM rows (rows 8, 9, 10) to Sheet1, matching the
# "added to BOM, is now complete" commit: a 10-pin female header,
# a 2-pin male header (with hyperlink), and the dual op-amp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 8: 10 pin female header
# ---------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A8").Value = "10 pin female header"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Value = "http://china.rs-online.com/web/p/pcb-sockets/7655745/"

$ws.Range("C8").Value = 2

$ws.Range("D8").Value = "J1, J2"

# ---------------------------------------------------------------
# Row 9: 2 pin male header (hyperlinked)
# ---------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Value = "2 pin male header"

$ws.Range("B4").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null

$h = $ws.Hyperlinks.Add($ws.Range("B9"), "http://china.rs-online.com//web/p/pcb-headers/2518086/", "", "", "http://china.rs-online.com/web/p/pcb-headers/2518086/")
$ws.Range("B9").Value = "http://china.rs-online.com//web/p/pcb-headers/2518086/"

# Hyperlinks.Add() resets the cell style; re-apply the wrapped hyperlink format.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null

$ws.Range("C9").Value = 2

$ws.Range("D9").Value = "J3, J4"

$ws.Rows.Item(9).RowHeight = 30

# ---------------------------------------------------------------
# Row 10: dual opamp, TL072CD
# ---------------------------------------------------------------
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = "dual opamp, TL072CD"

$ws.Range("B10").Value = "http://china.rs-online.com/web/p/operational-amplifiers/0857907/"

$ws.Range("C10").Value = 1

$ws.Range("D10").Value = "U1"

# ---------------------------------------------------------------
# Widen column D so the new "Place Labels" text fits.
# ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 76.5703125

# ---------------------------------------------------------------
# Leave the cursor where the user ended up after the edit.
# ---------------------------------------------------------------
$ws.Range("D12").Select() | Out-Null
